# Actualizacion de tareas - nueva tarea
#
# Adds a new task row to the "Hoja1" task list:
#   - C16 gets the status "en proceso" (reuses existing shared string).
#   - A30 gets the new task description "Validacion en creacion de cuota,
#     no muestra los mensajes de error" (previously an empty/placeholder cell).
#   - The sheet view's selection moves from C16 to B30, scrolled so row 13
#     is at the top of the viewport.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New status for the "Validacion de cuit" task row -> "en proceso"
$ws.Range("C16").Value = "en proceso"

# New task row text (row 30 was a blank placeholder row)
$ws.Range("A30").Value = "Validacion en creacion de cuota, no muestra los mensajes de error"
$ws.Range("A30").Font.Underline = $false

# Update the view: selection moves to B30, viewport scrolled to show row 13 at top
$ws.Range("B30").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
